$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '285.29'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.32%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.02%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.920'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.10%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06575'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.55%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.232'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '4.34%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.349'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '11.45%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9151'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '4.16%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1566'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.28%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06508'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '28.50%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07671'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.30%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02944'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.63%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08971'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.39%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001596'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.91%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006543'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.27%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006038'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.62%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.494'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.68%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.388'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.53%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.04%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3185'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.56%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1347'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.56%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.968'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.64%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1519'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '10.05%'
$ws.Range("B24").Value = 'CoinExToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04445'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.24%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001183'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.69%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004354'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '2.01%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001178'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-1.97%'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0001632'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '-15.82%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04156'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.03%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006852'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.54%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1414'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '20.49%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002036'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.97%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01180'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.44%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005542'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '6.93%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-2.63%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01846'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-7.80%'
